$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "29.220.88"
$ws.Range("E2").Value = "  -2.54%  "

# Row 3
$ws.Range("D3").Value = "1.853.32"
$ws.Range("E3").Value = "  -1.39%  "

# Row 4
$c = $ws.Range("D4")
$c.Value = "'0.9997"
$c.ClearFormats()
$ws.Range("E4").Value = "  -0.35%  "

# Row 5
$c = $ws.Range("D5")
$c.Value = "'0.6978"
$c.ClearFormats()

# Row 6
$c = $ws.Range("D6")
$c.Value = "'238.69"
$c.ClearFormats()
$ws.Range("E6").Value = "  -1.86%  "

# Row 7
$c = $ws.Range("D7")
$c.Value = "'0.9999"
$c.ClearFormats()
$ws.Range("E7").Value = "  -0.29%  "

# Row 8
$c = $ws.Range("D8")
$c.Value = "'0.3071"
$c.ClearFormats()
$ws.Range("E8").Value = "  -2.48%  "

# Row 9
$c = $ws.Range("D9")
$c.Value = "'0.07609"
$c.ClearFormats()
$ws.Range("E9").Value = "  +5.55%  "

# Row 10
$c = $ws.Range("D10")
$c.Value = "'23.63"
$c.ClearFormats()
$ws.Range("E10").Value = "  -4.21%  "

# Row 11
$c = $ws.Range("D11")
$c.Value = "'0.08088"
$c.ClearFormats()
$ws.Range("E11").Value = "  -3.25%  "

# Row 12
$ws.Range("B12").Value = "WrappedEther"
$ws.Range("C12").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D12").Value = "1.866.13"
$ws.Range("E12").Value = "  -1.46%  "

# Row 13
$ws.Range("B13").Value = "Polygon"
$ws.Range("C13").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$c = $ws.Range("D13")
$c.Value = "'0.7256"
$c.ClearFormats()
$ws.Range("E13").Value = "  -3.41%  "

# Row 14
$c = $ws.Range("D14")
$c.Value = "'5.188"
$c.ClearFormats()
$ws.Range("E14").Value = "  -4.16%  "

# Row 15
$c = $ws.Range("D15")
$c.Value = "'89.11"
$c.ClearFormats()

# Row 16
$ws.Range("D16").Value = "29.239.33"
$ws.Range("E16").Value = "  -2.51%  "

# Row 17
$c = $ws.Range("D17")
$c.Value = "'5.874"
$c.ClearFormats()
$ws.Range("E17").Value = "  -3.86%  "

# Row 18
$c = $ws.Range("D18")
$c.Value = "'241.76"
$c.ClearFormats()
$ws.Range("E18").Value = "  -2.88%  "

# Row 19
$c = $ws.Range("D19")
$c.Value = "'0.000007723"
$c.ClearFormats()
$ws.Range("E19").Value = "  -1.63%  "

# Row 20
$c = $ws.Range("D20")
$c.Value = "'13.12"
$c.ClearFormats()
$ws.Range("E20").Value = "  -3.34%  "

# Row 21
$c = $ws.Range("D21")
$c.Value = "'0.9999"
$c.ClearFormats()
$ws.Range("E21").Value = "  -0.14%  "

# Row 22
$ws.Range("D22").Value = "2.095.88"
$ws.Range("E22").Value = "  -2.46%  "

# Row 23
$c = $ws.Range("D23")
$c.Value = "'0.9999"
$c.ClearFormats()
$ws.Range("E23").Value = "  -0.36%  "

# Row 24
$c = $ws.Range("D24")
$c.Value = "'7.621"
$c.ClearFormats()
$ws.Range("E24").Value = "  -5.08%  "

# Row 25
$c = $ws.Range("D25")
$c.Value = "'9.056"
$c.ClearFormats()
$ws.Range("E25").Value = "  -2.37%  "

# Row 26
$c = $ws.Range("D26")
$c.Value = "'161.83"
$c.ClearFormats()
$ws.Range("E26").Value = "  -1.99%  "

# Row 27
$c = $ws.Range("D27")
$c.Value = "'0.1463"
$c.ClearFormats()
$ws.Range("E27").Value = "  -5.43%  "

# Row 28
$c = $ws.Range("D28")
$c.Value = "'18.06"
$c.ClearFormats()
$ws.Range("E28").Value = "  -3.61%  "

# Row 29
$c = $ws.Range("D29")
$c.Value = "'1.935"
$c.ClearFormats()
$ws.Range("E29").Value = "  -4.96%  "

# Row 30
$c = $ws.Range("D30")
$c.Value = "'1.400"
$c.ClearFormats()
$ws.Range("E30").Value = "  -7.77%  "

# Row 31
$ws.Range("B31").Value = "Filecoin"
$ws.Range("C31").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$c = $ws.Range("D31")
$c.Value = "'4.442"
$c.ClearFormats()
$ws.Range("E31").Value = "  -3.65%  "

# Row 32
$ws.Range("B32").Value = "PancakeSwap"
$ws.Range("C32").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$c = $ws.Range("D32")
$c.Value = "'1.502"
$c.ClearFormats()
$ws.Range("E32").Value = "  -2.13%  "

# Row 33
$c = $ws.Range("D33")
$c.Value = "'4.045"
$c.ClearFormats()
$ws.Range("E33").Value = "  -5.54%  "

# Row 34
$c = $ws.Range("D34")
$c.Value = "'0.05257"
$c.ClearFormats()
$ws.Range("E34").Value = "  -1.35%  "

# Row 35
$c = $ws.Range("D35")
$c.Value = "'1.193"
$c.ClearFormats()
$ws.Range("E35").Value = "  -3.53%  "

# Row 36
$c = $ws.Range("D36")
$c.Value = "'0.7114"
$c.ClearFormats()
$ws.Range("E36").Value = "  -5.15%  "

# Row 37
$c = $ws.Range("D37")
$c.Value = "'1.000"
$c.ClearFormats()
$ws.Range("E37").Value = "  -0.11%  "

# Row 38
$ws.Range("E38").Value = "  -1.52%  "

# Row 39
$c = $ws.Range("D39")
$c.Value = "'0.01863"
$c.ClearFormats()
$ws.Range("E39").Value = "  -5.28%  "

# Row 40
$c = $ws.Range("D40")
$c.Value = "'2.678"
$c.ClearFormats()
$ws.Range("E40").Value = "  -3.02%  "

# Row 41
$c = $ws.Range("D41")
$c.Value = "'0.9253"
$c.ClearFormats()
$ws.Range("E41").Value = "  +7.75%  "

# Row 42
$ws.Range("B42").Value = "FraxShare"
$ws.Range("C42").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$c = $ws.Range("D42")
$c.Value = "'5.949"
$c.ClearFormats()
$ws.Range("E42").Value = "  -3.01%  "

# Row 43
$ws.Range("B43").Value = "TheSandbox"
$ws.Range("C43").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$c = $ws.Range("D43")
$c.Value = "'0.4301"
$c.ClearFormats()
$ws.Range("E43").Value = "  -5.57%  "

# Row 44
$ws.Range("B44").Value = "Aave"
$ws.Range("C44").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$c = $ws.Range("D44")
$c.Value = "'69.71"
$c.ClearFormats()
$ws.Range("E44").Value = "  -3.65%  "

# Row 45
$ws.Range("B45").Value = "Maker"
$ws.Range("C45").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D45").Value = "1.046.00"
$ws.Range("E45").Value = "  -5.64%  "

# Row 46
$c = $ws.Range("D46")
$c.Value = "'0.9997"
$c.ClearFormats()
$ws.Range("E46").Value = "  -0.37%  "

# Row 47
$c = $ws.Range("D47")
$c.Value = "'102.33"
$c.ClearFormats()
$ws.Range("E47").Value = "  -1.90%  "

# Row 48
$c = $ws.Range("D48")
$c.Value = "'7.243"
$c.ClearFormats()
$ws.Range("E48").Value = "  -4.94%  "

# Row 49
$c = $ws.Range("D49")
$c.Value = "'1.741"
$c.ClearFormats()
$ws.Range("E49").Value = "  -6.20%  "

# Row 50
$c = $ws.Range("D50")
$c.Value = "'9.261"
$c.ClearFormats()
$ws.Range("E50").Value = "  -2.65%  "

# Row 51
$ws.Range("D51").Value = "1.991.91"
$ws.Range("E51").Value = "  -2.36%  "
